# Applies the "On 21 November 2016" update described in the commit:
#  - meta sheet: new note about picking which plants to keep
#  - vars sheet: new "Keep!" table listing the chosen varieties + totals

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Sheet "meta": add the dated note (rows 28-29)
# ---------------------------------------------------------------------
$meta = $wb.Worksheets.Item("meta")

$meta.Range("A28").Value = "On 21 November 2016:"
$meta.Range("A28").Font.Bold = $true

$meta.Range("B29").Value = "Lizzie updated this file to pick which plants to keep. She aimed for a mix of colors and phenophases from plants where at least some individuals flowered! (SeeVitisExpReps.xlsx)"

$meta.Activate() | Out-Null
$meta.Range("H35").Select() | Out-Null

# ---------------------------------------------------------------------
# Sheet "vars": add the "Keep!" summary table (rows 47-57)
# ---------------------------------------------------------------------
$vars = $wb.Worksheets.Item("vars")

$vars.Range("M47").Value = "Keep! For possible future experiments (21 Nov 2016):"
$vars.Range("M47").Font.Bold = $true

$vars.Range("R48").Value = "Reps total"

$vars.Range("M49").Value = "Gamay Noir"
$vars.Range("P49").Value = "early"
$vars.Range("Q49").Value = "red"
$vars.Range("R49").Value = 8

$vars.Range("M50").Value = "Cabernet Sauvignon"
$vars.Range("N50").Value = "Yes"
$vars.Range("P50").Value = "mid-late"
$vars.Range("Q50").Value = "red"
$vars.Range("R50").Value = 8

$vars.Range("M51").Value = "Syrah"
$vars.Range("N51").Value = "Yes"
$vars.Range("P51").Value = "mid"
$vars.Range("Q51").Value = "red"
$vars.Range("R51").Value = 7

$vars.Range("M52").Value = "Durif"
$vars.Range("P52").Value = "late"
$vars.Range("Q52").Value = "red"
$vars.Range("R52").Formula = "=7+6"

$vars.Range("M53").Value = "Tempranillo/Valdepenas"
$vars.Range("P53").Value = "early"
$vars.Range("Q53").Value = "red"
$vars.Range("R53").Value = 12

$vars.Range("M54").Value = "Verdelho"
$vars.Range("P54").Value = "very early"
$vars.Range("Q54").Value = "white"
$vars.Range("R54").Value = 6

$vars.Range("M55").Value = "Marsanne"
$vars.Range("P55").Value = "mid"
$vars.Range("Q55").Value = "white"
$vars.Range("R55").Value = 6

$vars.Range("R57").Formula = "=SUM(R49:R55)"

$vars.Activate() | Out-Null
$vars.Range("R58").Select() | Out-Null
